$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
  'D2' = '25.992.13'
  'E2' = '  +1.26%  '
  'D3' = '1.758.37'
  'E3' = '  +0.95%  '
  'D4' = '1.000'
  'E4' = '  -0.45%  '
  'D5' = '237.14'
  'E5' = '  -0.30%  '
  'D6' = '1.000'
  'E6' = '  -0.38%  '
  'D7' = '0.5222'
  'E7' = '  +3.88%  '
  'B8' = 'Cardano'
  'C8' = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
  'D8' = '0.2706'
  'E8' = '  +3.58%  '
  'B9' = 'Dogecoin'
  'C9' = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
  'D9' = '0.06200'
  'E9' = '  +1.19%  '
  'B10' = 'WrappedEther'
  'C10' = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
  'D10' = '1.769.24'
  'E10' = '  +1.32%  '
  'B11' = 'TRON'
  'C11' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
  'D11' = '0.07022'
  'E11' = '  +1.09%  '
  'B12' = 'Solana'
  'C12' = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
  'D12' = '15.71'
  'E12' = '  +4.00%  '
  'B13' = 'Polygon'
  'C13' = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
  'D13' = '0.6568'
  'E13' = '  +11.27%  '
  'B14' = 'Polkadot'
  'C14' = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
  'D14' = '4.483'
  'E14' = '  +0.29%  '
  'B15' = 'Litecoin'
  'C15' = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
  'D15' = '78.17'
  'E15' = '  +2.31%  '
  'B16' = 'BinanceUSD'
  'C16' = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
  'D16' = '1.001'
  'E16' = '  -0.31%  '
  'B17' = 'Dai'
  'C17' = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
  'D17' = '0.9996'
  'E17' = '  -0.31%  '
  'B18' = 'WrappedBTC'
  'C18' = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
  'D18' = '25.992.59'
  'E18' = '  +1.02%  '
  'B19' = 'Avalanche'
  'C19' = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
  'D19' = '11.68'
  'E19' = '  +0.99%  '
  'B20' = 'ShibaInu'
  'C20' = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
  'D20' = '0.000006699'
  'E20' = '  -0.57%  '
  'B21' = 'WrappedliquidstakedEther2.0'
  'C21' = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
  'D21' = '1.983.79'
  'E21' = '  +0.67%  '
  'B22' = 'Uniswap'
  'C22' = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
  'D22' = '4.093'
  'E22' = '  +1.42%  '
  'B23' = 'Cosmos'
  'C23' = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
  'D23' = '8.407'
  'E23' = '  +4.19%  '
  'B24' = 'Chainlink'
  'C24' = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
  'D24' = '5.176'
  'E24' = '  +1.90%  '
  'B25' = 'Monero'
  'C25' = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
  'D25' = '137.25'
  'E25' = '  -0.67%  '
  'B26' = 'Toncoin'
  'C26' = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
  'D26' = '1.482'
  'E26' = '  -3.36%  '
  'B27' = 'EthereumClassic'
  'C27' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
  'D27' = '15.14'
  'E27' = '  +1.55%  '
  'D28' = '1.828'
  'E28' = '  +1.32%  '
  'B29' = 'BitcoinCash'
  'C29' = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
  'D29' = '102.53'
  'E29' = '  -0.39%  '
  'B30' = 'Stellar'
  'C30' = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
  'D30' = '0.08408'
  'E30' = '  +3.79%  '
  'B31' = 'InternetComputer(DFINITY)'
  'C31' = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
  'D31' = '3.698'
  'E31' = '  -1.60%  '
  'B32' = 'Filecoin'
  'C32' = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
  'D32' = '3.417'
  'E32' = '  -0.51%  '
  'B33' = 'Hedera'
  'C33' = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
  'D33' = '0.04419'
  'E33' = '  -1.39%  '
  'B34' = 'HuobiToken'
  'C34' = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
  'D34' = '2.650'
  'E34' = '  +0.52%  '
  'B35' = 'ARBITRUM'
  'C35' = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
  'D35' = '0.9966'
  'E35' = '  +3.12%  '
  'B36' = 'ImmutableX'
  'C36' = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
  'D36' = '0.6078'
  'E36' = '  +1.10%  '
  'B37' = 'MXToken'
  'C37' = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
  'D37' = '2.736'
  'E37' = '  +3.16%  '
  'B38' = 'VeChain'
  'C38' = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
  'D38' = '0.01571'
  'E38' = '  +1.91%  '
  'B39' = 'RenderToken'
  'C39' = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
  'D39' = '1.949'
  'E39' = '  +2.56%  '
  'B40' = 'PaxDollar'
  'C40' = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
  'D40' = '1.001'
  'E40' = '  -0.12%  '
  'B41' = 'Quant'
  'C41' = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
  'D41' = '102.92'
  'E41' = '  -1.09%  '
  'B42' = 'TheSandbox'
  'C42' = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
  'D42' = '0.3881'
  'E42' = '  +2.95%  '
  'B43' = 'TrustWalletToken'
  'C43' = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
  'D43' = '0.7525'
  'E43' = '  +3.81%  '
  'B44' = 'FraxShare'
  'C44' = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
  'D44' = '4.933'
  'E44' = '  -3.64%  '
  'B45' = 'Cronos'
  'C45' = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
  'D45' = '0.05496'
  'E45' = '  +3.13%  '
  'B46' = 'Algorand'
  'C46' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
  'D46' = '0.1118'
  'E46' = '  +1.48%  '
  'B47' = 'Aptos'
  'C47' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
  'D47' = '6.101'
  'E47' = '  +4.36%  '
  'B48' = 'Elrond'
  'C48' = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
  'D48' = '30.14'
  'E48' = '  +0.90%  '
  'B49' = 'Aave'
  'C49' = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
  'D49' = '52.64'
  'E49' = '  +0.91%  '
  'B50' = 'USDD'
  'C50' = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
  'D50' = '1.002'
  'E50' = '  -0.18%  '
  'B51' = 'Decentraland'
  'C51' = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
  'D51' = '0.3432'
  'E51' = '  +0.34%  '
}

foreach ($addr in $updates.Keys) {
  $cell = $ws.Range($addr)
  $cell.NumberFormat = "@"
  $cell.Value = $updates[$addr]
  $cell.ClearFormats()
}